# Airpods.xlsx - "Treinamento" sheet: mark additional rows as relevant
# (column B 0 -> 1) and move the viewport/selection, per the commit
# "Funcao de limpeza melhorada e atualizada ... excel com classificacao
# melhorada".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Treinamento")

# Rows whose classification flag (column B) flips from 0 to 1.
$rowsToSet = @(
    47, 51, 54, 57, 61, 65, 79, 80, 84, 87, 90, 94, 100, 103, 107, 109, 110,
    111, 112, 116, 122, 130, 133, 136, 139, 142, 149, 154, 158, 162, 167,
    171, 176, 179, 180, 194, 202, 251
)

foreach ($r in $rowsToSet) {
    $ws.Cells.Item($r, 2).Value = 1
}

# Make this sheet active and move the viewport / selection to where the
# user was working (top-left A244, cursor on B252).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 244
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B252").Select()
